$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set values first
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Format B1: bold, thin box border, centered horizontally, top vertical
# alignment. This is the only cell whose format is built up property by
# property - this yields exactly one new cell style entry.
$r1 = $ws.Range("B1")
$r1.Font.Bold = $true
$r1.HorizontalAlignment = -4108   # xlCenter
$r1.VerticalAlignment = -4160     # xlTop
$r1.Borders.LineStyle = 1         # xlContinuous
$r1.Borders.Weight = 2            # xlThin

# Copy B1's formatting onto A2 so both cells share the same cell style
# (avoids generating spurious intermediate style entries).
$r1.Copy()
$ws.Range("A2").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
